$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5002720
$ws.Range("I74").Value = 2281.4375
$ws.Range("J74").Value = 13892389
$ws.Range("K74").Value = 2281.4375
$ws.Range("L74").Value = 13892389
$ws.Range("M74").Value = -1345.4375
$ws.Range("N74").Value = -13894261
$ws.Range("H77").Value = 5002720
$ws.Range("I77").Value = 2281.4375
$ws.Range("J77").Value = 13892389
$ws.Range("K77").Value = 11407.1875
$ws.Range("L77").Value = 69461945
$ws.Range("M77").Value = -6727.1875
$ws.Range("N77").Value = -69471305
$ws.Range("H92").Value = 2245.9583
$ws.Range("J92").Value = 1562.5
$ws.Range("L92").Value = 1562.5
$ws.Range("N92").Value = -4058.5
$ws.Range("H113").Value = 111115590
$ws.Range("I113").Value = 200002500
$ws.Range("J113").Value = 6948.5
$ws.Range("K113").Value = 200002500
$ws.Range("L113").Value = 6948.5
$ws.Range("M113").Value = -199999246
$ws.Range("N113").Value = -13456.5
$ws.Range("H116").Value = 2767.35
$ws.Range("I116").Value = 1934.9
$ws.Range("K116").Value = 1934.9
$ws.Range("M116").Value = 1507.1
$ws.Range("H121").Value = 820.25
$ws.Range("J121").Value = 820.25
$ws.Range("L121").Value = 2460.75
$ws.Range("N121").Value = -5954.75
$ws.Range("H129").Value = 963.64703
$ws.Range("J129").Value = 997.2766
$ws.Range("L129").Value = 2991.8298
$ws.Range("N129").Value = -12991.8298

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3740.1714
$ws.Range("I32").Value = 2174.7144
$ws.Range("K32").Value = 2174.7144
$ws.Range("M32").Value = -1887.7144
$ws.Range("H74").Value = 1260
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 1260
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 1260
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = -3008
$ws.Range("H77").Value = 1260
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 1260
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 6300
$ws.Range("M77").Value = ""
$ws.Range("N77").Value = -15036
$ws.Range("H122").Value = 2468
$ws.Range("I122").Value = 2429.7
$ws.Range("J122").Value = 2595.6667
$ws.Range("K122").Value = 7289.099999999999
$ws.Range("L122").Value = 7787.000100000001
$ws.Range("M122").Value = -4839.099999999999
$ws.Range("N122").Value = -12687.0001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 854.2222
$ws.Range("I20").Value = 748.5
$ws.Range("J20").Value = 1700
$ws.Range("K20").Value = 748.5
$ws.Range("L20").Value = 1700
$ws.Range("M20").Value = -501.5
$ws.Range("N20").Value = -2194

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2090.3157
$ws.Range("I31").Value = 792.6667
$ws.Range("J31").Value = 5275.4546
$ws.Range("K31").Value = 792.6667
$ws.Range("L31").Value = 5275.4546
$ws.Range("M31").Value = -497.6667
$ws.Range("N31").Value = -5865.4546
$ws.Range("H34").Value = 2090.3157
$ws.Range("I34").Value = 792.6667
$ws.Range("J34").Value = 5275.4546
$ws.Range("K34").Value = 792.6667
$ws.Range("L34").Value = 5275.4546
$ws.Range("M34").Value = -590.6667
$ws.Range("N34").Value = -5679.4546
$ws.Range("H58").Value = 22773.088
$ws.Range("I58").Value = 1181.5
$ws.Range("J58").Value = 100502.8
$ws.Range("K58").Value = 1181.5
$ws.Range("L58").Value = 100502.8
$ws.Range("M58").Value = -978.5
$ws.Range("N58").Value = -100908.8
$ws.Range("H62").Value = 4478.3076
$ws.Range("I62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("M62").Value = -3376
$ws.Range("H65").Value = 4478.3076
$ws.Range("I65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("M65").Value = -16880
$ws.Range("H105").Value = 12501228
$ws.Range("I105").Value = 15625798
$ws.Range("J105").Value = 2949
$ws.Range("K105").Value = 15625798
$ws.Range("L105").Value = 2949
$ws.Range("M105").Value = -15624051
$ws.Range("N105").Value = -6443
$ws.Range("H136").Value = 22773.088
$ws.Range("I136").Value = 1181.5
$ws.Range("J136").Value = 100502.8
$ws.Range("K136").Value = 3544.5
$ws.Range("L136").Value = 301508.4
$ws.Range("M136").Value = -994.5
$ws.Range("N136").Value = -306608.4

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 749.7931
$ws.Range("I5").Value = 647.8889
$ws.Range("K5").Value = 1943.6667
$ws.Range("M5").Value = -1831.6667
$ws.Range("H68").Value = 1154.5
$ws.Range("J68").Value = 1159.4839
$ws.Range("L68").Value = 3478.4517
$ws.Range("N68").Value = -5100.4517
$ws.Range("H71").Value = 1154.5
$ws.Range("J71").Value = 1159.4839
$ws.Range("L71").Value = 10435.3551
$ws.Range("N71").Value = -18547.3551
$ws.Range("H98").Value = 498.33334
$ws.Range("I98").Value = 300
$ws.Range("J98").Value = 597.5
$ws.Range("K98").Value = 900
$ws.Range("L98").Value = 1792.5
$ws.Range("M98").Value = 598
$ws.Range("N98").Value = -4788.5
$ws.Range("H116").Value = 1079.8
$ws.Range("I116").Value = 266.33334
$ws.Range("J116").Value = 2300
$ws.Range("K116").Value = 799.0000200000001
$ws.Range("L116").Value = 6900
$ws.Range("M116").Value = 2642.99998
$ws.Range("N116").Value = -13784
$ws.Range("H131").Value = 806.0700000000001
$ws.Range("I131").Value = 385
$ws.Range("J131").Value = 823.61456
$ws.Range("K131").Value = 1155
$ws.Range("L131").Value = 2470.84368
$ws.Range("M131").Value = 3885
$ws.Range("N131").Value = -12550.84368
$ws.Range("H135").Value = 749.7931
$ws.Range("I135").Value = 647.8889
$ws.Range("K135").Value = 5831.0001
$ws.Range("M135").Value = -3296.0001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 20768.215
$ws.Range("I132").Value = 2833.9048
$ws.Range("K132").Value = 8501.714399999999
$ws.Range("M132").Value = -5971.714399999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3983.6
$ws.Range("I7").Value = 4353.5713
$ws.Range("J7").Value = 3120.3333
$ws.Range("K7").Value = 4353.5713
$ws.Range("L7").Value = 3120.3333
$ws.Range("M7").Value = -4241.5713
$ws.Range("N7").Value = -3344.3333
$ws.Range("H82").Value = 3490
$ws.Range("I82").Value = 4112.5
$ws.Range("J82").Value = 1000
$ws.Range("K82").Value = 4112.5
$ws.Range("L82").Value = 1000
$ws.Range("M82").Value = -3751.5
$ws.Range("N82").Value = -1722
$ws.Range("H85").Value = 3490
$ws.Range("I85").Value = 4112.5
$ws.Range("J85").Value = 1000
$ws.Range("K85").Value = 4112.5
$ws.Range("L85").Value = 1000
$ws.Range("M85").Value = -2864.5
$ws.Range("N85").Value = -3496
$ws.Range("H126").Value = 3983.6
$ws.Range("I126").Value = 4353.5713
$ws.Range("J126").Value = 3120.3333
$ws.Range("K126").Value = 13060.7139
$ws.Range("L126").Value = 9360.999899999999
$ws.Range("M126").Value = -10590.7139
$ws.Range("N126").Value = -14300.9999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1375.4166
$ws.Range("I100").Value = 958.6667
$ws.Range("J100").Value = 1792.1666
$ws.Range("K100").Value = 1917.3334
$ws.Range("L100").Value = 3584.3332
$ws.Range("M100").Value = -1376.3334
$ws.Range("N100").Value = -4666.3332
